$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / subtitle text (October -> November) ---
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Million Cubic Feet)"

# --- Insert a new row for the "November" monthly data in the "Year 2016" block ---
# (this pushes the old rows 53-60 down to 54-61)
$ws.Rows("53").Insert()

# Copy formatting from the row above (October, row 52) into the new row so that
# styles/borders match the rest of the monthly data rows.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new November data row.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 395
$ws.Range("C53").Value = 0.3
$ws.Range("D53").Value = 214
$ws.Range("E53").Value = 90
$ws.Range("F53").Value = 90

# --- Update "Rolling 12 Months Ending in October" -> "...November" header ---
# (this row was 57 before the insert, it is now row 58)
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# --- Update "Year to Date" figures (rows 54 now shifted to 55-57) ---
# Year 2014
$ws.Range("B55").Value = 1577
$ws.Range("C55").Value = 161
$ws.Range("D55").Value = 499
$ws.Range("E55").Value = 613
$ws.Range("F55").Value = 305

# Year 2015
$ws.Range("B56").Value = 1370
$ws.Range("C56").Value = 2
$ws.Range("D56").Value = 568
$ws.Range("E56").Value = 470
$ws.Range("F56").Value = 330

# Year 2016
$ws.Range("B57").Value = 3814
$ws.Range("C57").Value = 5
$ws.Range("D57").Value = 2152
$ws.Range("E57").Value = 875
$ws.Range("F57").Value = 783

# --- Update "Rolling 12 Months" figures (rows 59-60) ---
# 2015
$ws.Range("B59").Value = 1503
$ws.Range("C59").Value = 17
$ws.Range("D59").Value = 594
$ws.Range("E59").Value = 531
$ws.Range("F59").Value = 361

# 2016
$ws.Range("B60").Value = 3966
$ws.Range("C60").Value = 5
$ws.Range("D60").Value = 2227
$ws.Range("E60").Value = 919
$ws.Range("F60").Value = 815
